$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.878.71"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "3.454.10"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").Value = "'159.49"
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.453.28"
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").Value = "'0.574"
$ws.Range("E9").Value = "  -6.25%  "
$ws.Range("D10").Value = "'7.20"
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("D11").Value = "'0.122"
$ws.Range("E11").Value = "  -3.31%  "
$ws.Range("D12").Value = "'0.440"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("D13").Value = "4.048.92"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "'27.65"
$ws.Range("E15").Value = "  -4.18%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000173"
$ws.Range("E16").Value = "  -10.77%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "64.934.33"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "3.453.49"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "'6.21"
$ws.Range("E19").Value = "  -4.19%  "
$ws.Range("D20").Value = "'13.71"
$ws.Range("E20").Value = "  -5.11%  "
$ws.Range("D21").Value = "'376.79"
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("E22").Value = "  -3.84%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'72.22"
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("E25").Value = "  -3.52%  "
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").Value = "'9.93"
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -5.90%  "
$ws.Range("D31").Value = "'6.07"
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("E32").Value = "  -2.47%  "
$ws.Range("D33").Value = "'23.19"
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("D34").Value = "'7.00"
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("E35").Value = "  -4.23%  "
$ws.Range("D36").Value = "'161.27"
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("D37").Value = "'1.88"
$ws.Range("E37").Value = "  -2.61%  "
$ws.Range("D38").Value = "2.895.18"
$ws.Range("E38").Value = "  -4.06%  "
$ws.Range("D39").Value = "'0.0749"
$ws.Range("E39").Value = "  -4.58%  "
$ws.Range("D40").Value = "'26.32"
$ws.Range("E40").Value = "  -2.92%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'4.52"
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'43.04"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "'0.788"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").Value = "'26.19"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "'6.38"
$ws.Range("E45").Value = "  -6.71%  "
$ws.Range("D46").Value = "'0.0310"
$ws.Range("E46").Value = "  -3.87%  "
$ws.Range("E47").Value = "  +8.31%  "
$ws.Range("D48").Value = "'321.31"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  -3.19%  "
$ws.Range("D50").Value = "'6.47"
$ws.Range("E50").Value = "  -4.27%  "
$ws.Range("D51").Value = "'0.846"
$ws.Range("E51").Value = "  -4.44%  "
